$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
}

Set-TextCell 'D2' '66.827.97'
Set-TextCell 'E2' '  -1.30%  '
Set-TextCell 'D3' '3.221.26'
Set-TextCell 'E3' '  -2.62%  '
Set-TextCell 'D4' '1.00'
Set-TextCell 'E4' '  +0.02%  '
Set-TextCell 'D5' '581.79'
Set-TextCell 'E5' '  -3.12%  '
Set-TextCell 'D6' '141.17'
Set-TextCell 'E6' '  -12.50%  '
Set-TextCell 'D7' '0.998'
Set-TextCell 'D8' '3.214.11'
Set-TextCell 'E8' '  -2.63%  '
Set-TextCell 'D9' '0.525'
Set-TextCell 'E9' '  -7.81%  '
Set-TextCell 'D10' '0.162'
Set-TextCell 'E10' '  -10.53%  '
Set-TextCell 'D11' '6.37'
Set-TextCell 'E11' '  -2.52%  '
Set-TextCell 'D12' '0.479'
Set-TextCell 'E12' '  -8.34%  '
Set-TextCell 'D13' '0.0000234'
Set-TextCell 'E13' '  -7.79%  '
Set-TextCell 'D14' '35.88'
Set-TextCell 'E14' '  -12.65%  '
Set-TextCell 'D15' '3.726.67'
Set-TextCell 'E15' '  -3.36%  '
Set-TextCell 'D16' '66.861.17'
Set-TextCell 'E16' '  -1.44%  '
Set-TextCell 'D17' '3.207.49'
Set-TextCell 'E17' '  -3.48%  '
Set-TextCell 'E18' '  -3.88%  '
Set-TextCell 'D19' '505.39'
Set-TextCell 'E19' '  -7.43%  '
Set-TextCell 'D20' '6.80'
Set-TextCell 'E20' '  -10.55%  '
Set-TextCell 'D21' '14.24'
Set-TextCell 'E21' '  -9.98%  '
Set-TextCell 'D22' '0.715'
Set-TextCell 'E22' '  -9.53%  '
Set-TextCell 'D23' '7.37'
Set-TextCell 'E23' '  -10.35%  '
Set-TextCell 'D24' '81.74'
Set-TextCell 'E24' '  -8.15%  '
Set-TextCell 'D25' '12.84'
Set-TextCell 'E25' '  -9.09%  '
Set-TextCell 'E26' '  +0.47%  '
Set-TextCell 'D27' '3.07'
Set-TextCell 'E27' '  -10.57%  '
Set-TextCell 'B28' 'ImmutableX'
Set-TextCell 'C28' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D28' '2.04'
Set-TextCell 'E28' '  -9.44%  '
Set-TextCell 'B29' 'EthereumClassic'
Set-TextCell 'C29' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D29' '27.82'
Set-TextCell 'E29' '  -8.98%  '
Set-TextCell 'D30' '7.51'
Set-TextCell 'E30' '  -6.39%  '
Set-TextCell 'E31' '  -1.29%  '
Set-TextCell 'D32' '2.50'
Set-TextCell 'E32' '  -4.86%  '
Set-TextCell 'D33' '1.00'
Set-TextCell 'E33' '  -0.45%  '
Set-TextCell 'D34' '6.09'
Set-TextCell 'E34' '  -15.44%  '
Set-TextCell 'D35' '496.17'
Set-TextCell 'E35' '  -14.54%  '
Set-TextCell 'D36' '54.07'
Set-TextCell 'E36' '  -2.05%  '
Set-TextCell 'D37' '5.30'
Set-TextCell 'E37' '  -12.77%  '
Set-TextCell 'D38' '0.0414'
Set-TextCell 'E38' '  -6.40%  '
Set-TextCell 'D39' '0.0814'
Set-TextCell 'E39' '  -9.39%  '
Set-TextCell 'D40' '8.50'
Set-TextCell 'E40' '  -12.73%  '
Set-TextCell 'D41' '0.119'
Set-TextCell 'E41' '  -12.14%  '
Set-TextCell 'D42' '2.828.63'
Set-TextCell 'E42' '  -6.35%  '
Set-TextCell 'D43' '2.54'
Set-TextCell 'E43' '  -9.54%  '
Set-TextCell 'B44' 'TheGraph'
Set-TextCell 'C44' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell 'D44' '0.252'
Set-TextCell 'E44' '  -8.20%  '
Set-TextCell 'B45' 'USDe'
Set-TextCell 'C45' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell 'D45' '1.00'
Set-TextCell 'E45' '  -0.10%  '
Set-TextCell 'D46' '25.05'
Set-TextCell 'E46' '  -11.79%  '
Set-TextCell 'B47' 'Monero'
Set-TextCell 'C47' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D47' '121.44'
Set-TextCell 'E47' '  -4.65%  '
Set-TextCell 'B48' 'Fetch.AI'
Set-TextCell 'C48' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 'D48' '2.03'
Set-TextCell 'E48' '  -8.30%  '
Set-TextCell 'D49' '0.0₃0529'
Set-TextCell 'E49' '  -13.77%  '
Set-TextCell 'D50' '0.109'
Set-TextCell 'E50' '  -8.64%  '
Set-TextCell 'D51' '2.13'
Set-TextCell 'E51' '  -18.70%  '
